# Apply updated Universalis market-price pricing figures to each job leve-profit sheet.
# Values recomputed by the scheduled pricing runner; only columns H-N (computed price/profit
# columns) change - the leve metadata in columns A-G is untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 4 (Leve Item ID 5470)
$ws.Range("H4").Value = 1350
$ws.Range("I4").Value = 1630
$ws.Range("J4").Value = 650
$ws.Range("K4").Value = 1630
$ws.Range("L4").Value = 650
$ws.Range("M4").Value = -1516
$ws.Range("N4").Value = -878

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 79854.57000000001
$ws.Range("J43").Value = 78304.92
$ws.Range("L43").Value = 78304.92
$ws.Range("N43").Value = -78442.92

# Row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 2843.1428
$ws.Range("J70").Value = 2957.4285
$ws.Range("L70").Value = 8872.2855
$ws.Range("N70").Value = -9412.2855

# Row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 2843.1428
$ws.Range("J73").Value = 2957.4285
$ws.Range("L73").Value = 8872.2855
$ws.Range("N73").Value = -10744.2855

# Row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 1101.6
$ws.Range("I92").Value = 1101.6
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1101.6
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 146.4000000000001
$ws.Range("N92").ClearContents()

# Row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 1433.3334
$ws.Range("I111").Value = 1433.3334
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 4300.0002
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -1233.0002
$ws.Range("N111").ClearContents()

# Row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 1782
$ws.Range("I113").Value = 1740
$ws.Range("J113").Value = 1803
$ws.Range("K113").Value = 1740
$ws.Range("L113").Value = 1803
$ws.Range("M113").Value = 1514
$ws.Range("N113").Value = -8311

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 4351.303
$ws.Range("I137").Value = 1075.25
$ws.Range("J137").Value = 4803.1724
$ws.Range("K137").Value = 3225.75
$ws.Range("L137").Value = 14409.5172
$ws.Range("M137").Value = -675.75
$ws.Range("N137").Value = -19509.5172

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 1690.5952
$ws.Range("I138").Value = 820.2778
$ws.Range("J138").Value = 1927.9546
$ws.Range("K138").Value = 2460.8334
$ws.Range("L138").Value = 5783.8638
$ws.Range("M138").Value = 2679.1666
$ws.Range("N138").Value = -16063.8638

# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 2684.2222
$ws.Range("I141").Value = 2221.7058
$ws.Range("K141").Value = 6665.117400000001
$ws.Range("M141").Value = -1485.117400000001

$ws = $wb.Worksheets.Item("ARM")

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 2197.0527
$ws.Range("I45").Value = 1947.3846
$ws.Range("K45").Value = 1947.3846
$ws.Range("M45").Value = -1570.3846

$ws = $wb.Worksheets.Item("BSM")

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 501948.72
$ws.Range("I86").Value = 1978.2
$ws.Range("J86").Value = 1751875
$ws.Range("K86").Value = 1978.2
$ws.Range("L86").Value = 1751875
$ws.Range("M86").Value = -855.2
$ws.Range("N86").Value = -1754121

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 501948.72
$ws.Range("I89").Value = 1978.2
$ws.Range("J89").Value = 1751875
$ws.Range("K89").Value = 9891
$ws.Range("L89").Value = 8759375
$ws.Range("M89").Value = -4275
$ws.Range("N89").Value = -8770607

# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 1106
$ws.Range("I99").Value = 868.3333
$ws.Range("J99").Value = 1462.5
$ws.Range("K99").Value = 868.3333
$ws.Range("L99").Value = 1462.5
$ws.Range("M99").Value = 629.6667
$ws.Range("N99").Value = -4458.5

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 34049.22
$ws.Range("I134").Value = 1748.5
$ws.Range("K134").Value = 5245.5
$ws.Range("M134").Value = -2710.5

$ws = $wb.Worksheets.Item("CRP")

# Row 94 (Leve Item ID 32934)
$ws.Range("H94").Value = 1509.909
$ws.Range("I94").Value = 999.25
$ws.Range("J94").Value = 1801.7142
$ws.Range("K94").Value = 999.25
$ws.Range("L94").Value = 1801.7142
$ws.Range("M94").Value = -548.25
$ws.Range("N94").Value = -2703.7142

# Row 120 (Leve Item ID 27230)
$ws.Range("H120").Value = 30244.5
$ws.Range("J120").Value = 30244.5
$ws.Range("L120").Value = 30244.5
$ws.Range("N120").Value = -37502.5

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2195.1304
$ws.Range("I132").Value = 1450
$ws.Range("K132").Value = 4350
$ws.Range("M132").Value = -1820

$ws = $wb.Worksheets.Item("CUL")

# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 6325.6113
$ws.Range("I5").Value = 658.44446
$ws.Range("J5").Value = 11992.777
$ws.Range("K5").Value = 1975.33338
$ws.Range("L5").Value = 35978.331
$ws.Range("M5").Value = -1863.33338
$ws.Range("N5").Value = -36202.331

# Row 117 (Leve Item ID 27870)
$ws.Range("H117").Value = 850.63635
$ws.Range("I117").Value = 445.8
$ws.Range("J117").Value = 1188
$ws.Range("K117").Value = 1337.4
$ws.Range("L117").Value = 3564
$ws.Range("M117").Value = 2104.6
$ws.Range("N117").Value = -10448

# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 423.10526
$ws.Range("I122").Value = 295.7143
$ws.Range("J122").Value = 779.8
$ws.Range("K122").Value = 2661.4287
$ws.Range("L122").Value = 7018.2
$ws.Range("M122").Value = -211.4286999999999
$ws.Range("N122").Value = -11918.2

# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1038.8889
$ws.Range("I132").Value = 912.5
$ws.Range("J132").Value = 1140
$ws.Range("K132").Value = 8212.5
$ws.Range("L132").Value = 10260
$ws.Range("M132").Value = -5682.5
$ws.Range("N132").Value = -15320

# Row 133 (Leve Item ID 44073)
$ws.Range("H133").Value = 3293
$ws.Range("I133").Value = 1380
$ws.Range("J133").Value = 3930.6667
$ws.Range("K133").Value = 4140
$ws.Range("L133").Value = 11792.0001
$ws.Range("M133").Value = 920
$ws.Range("N133").Value = -21912.0001

# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 6325.6113
$ws.Range("I135").Value = 658.44446
$ws.Range("J135").Value = 11992.777
$ws.Range("K135").Value = 5926.00014
$ws.Range("L135").Value = 107934.993
$ws.Range("M135").Value = -3391.00014
$ws.Range("N135").Value = -113004.993

$ws = $wb.Worksheets.Item("GSM")

# Row 39 (Leve Item ID 18264)
$ws.Range("H39").Value = 27200
$ws.Range("J39").Value = 27200
$ws.Range("L39").Value = 27200
$ws.Range("N39").Value = -28264

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 2008.7693
$ws.Range("I126").Value = 1533.3334
$ws.Range("J126").Value = 2151.4
$ws.Range("K126").Value = 4600.0002
$ws.Range("L126").Value = 6454.200000000001
$ws.Range("M126").Value = -2130.0002
$ws.Range("N126").Value = -11394.2

$ws = $wb.Worksheets.Item("LTW")

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 1021.35
$ws.Range("I22").Value = 803.3333
$ws.Range("J22").Value = 1199.7273
$ws.Range("K22").Value = 803.3333
$ws.Range("L22").Value = 1199.7273
$ws.Range("M22").Value = -508.3333
$ws.Range("N22").Value = -1789.7273

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 1021.35
$ws.Range("I27").Value = 803.3333
$ws.Range("J27").Value = 1199.7273
$ws.Range("K27").Value = 803.3333
$ws.Range("L27").Value = 1199.7273
$ws.Range("M27").Value = -696.3333
$ws.Range("N27").Value = -1413.7273

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 1550.25
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1550.25
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1550.25
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1926.25

# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 933.9524
$ws.Range("I55").Value = 449.15384
$ws.Range("J55").Value = 1721.75
$ws.Range("K55").Value = 449.15384
$ws.Range("L55").Value = 1721.75
$ws.Range("M55").Value = -276.15384
$ws.Range("N55").Value = -2067.75

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 2077.3076
$ws.Range("I122").Value = 2460
$ws.Range("J122").Value = 1838.125
$ws.Range("K122").Value = 7380
$ws.Range("L122").Value = 5514.375
$ws.Range("M122").Value = -4930
$ws.Range("N122").Value = -10414.375

$ws = $wb.Worksheets.Item("WVR")

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 909.5217
$ws.Range("I107").Value = 718.2308
$ws.Range("J107").Value = 1158.2
$ws.Range("K107").Value = 2154.6924
$ws.Range("L107").Value = 3474.6
$ws.Range("M107").Value = -234.6923999999999
$ws.Range("N107").Value = -7314.6

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 1229.3125
$ws.Range("I126").Value = 1066.9412
$ws.Range("J126").Value = 1413.3334
$ws.Range("K126").Value = 3200.8236
$ws.Range("L126").Value = 4240.0002
$ws.Range("M126").Value = -730.8235999999997
$ws.Range("N126").Value = -9180.0002
